# Automated "availability check" run: refreshes the timestamp on the last
# existing batch of rows (982:995) and appends one more 14-row batch
# (996:1009) with a newer timestamp, mirroring the recurring block pattern
# already present in the sheet (one row per monitored service).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Touch the previous batch's timestamp (982:995) ---------------------
$refreshedTimestamp = 44233.15125248842
for ($r = 982; $r -le 995; $r++) {
  $ws.Range("D$r").Value = $refreshedTimestamp
}

# --- 2) Append a new batch of 14 rows (996:1009) ----------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")

$displayUrls = @(
  "https://www.dataintelligence-group.com/",
  "https://serviciodashboard.azurewebsites.net/",
  "https://powerbi.microsoft.com/es-es/",
  "https://www.dropbox.com/",
  "https://dataintelligence.store/",
  "https://app-data-i.users.earthengine.app/",
  "https://odooutil.azurewebsites.net/",
  "https://filtradordashboard.azurewebsites.net/",
  "https://ide.dataintelligence-group.com/mapstore/#/",
  "https://ide.dataintelligence-group.com/geoserver/web/?0",
  "https://ide.dataintelligence-group.com/",
  "https://rpubs.com/dataintelligence/",
  "https://github.com/Sud-Austral/",
  "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$linkAddresses = @(
  "https://www.dataintelligence-group.com/",
  "https://serviciodashboard.azurewebsites.net/",
  "https://powerbi.microsoft.com/es-es/",
  "https://www.dropbox.com/",
  "https://dataintelligence.store/",
  "https://app-data-i.users.earthengine.app/",
  "https://odooutil.azurewebsites.net/",
  "https://filtradordashboard.azurewebsites.net/",
  "https://ide.dataintelligence-group.com/mapstore/",
  "https://ide.dataintelligence-group.com/geoserver/web/?0",
  "https://ide.dataintelligence-group.com/",
  "https://rpubs.com/dataintelligence/",
  "https://github.com/Sud-Austral/",
  "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$linkSubAddresses = @("","","","","","","","","/","","","","","")

$availability = "Disponible"
$newTimestamp = 44233.1725618268
$startRow = 996

for ($i = 0; $i -lt $names.Count; $i++) {
  $row = $startRow + $i

  $ws.Range("A$row").Value = $names[$i]
  $ws.Range("B$row").Value = $displayUrls[$i]
  $ws.Range("C$row").Value = $availability

  $ws.Range("D$row").Value = $newTimestamp
  $ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

  if ($linkSubAddresses[$i] -ne "") {
    $ws.Hyperlinks.Add($ws.Range("B$row"), $linkAddresses[$i], $linkSubAddresses[$i])
  } else {
    $ws.Hyperlinks.Add($ws.Range("B$row"), $linkAddresses[$i])
  }

  # Reapply the sheet's standard hyperlink cell style (Add() leaves its own
  # default formatting instead of reusing the shared "Hyperlink" style).
  $ws.Range("B$row").Style = $ws.Range("B2").Style
}
